$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the now-empty C62 cell entirely (was an empty styled cell, removed in target)
$ws.Range("C62").Clear()

# Update duration/step-count cells in column C from time-range strings to plain numbers
$ws.Range("C38").Value = 4
$ws.Range("C39").Value = 6
$ws.Range("C40").Value = 8
$ws.Range("C41").Value = 13
$ws.Range("C42").Value = 7
$ws.Range("C43").Value = 22
$ws.Range("C47").Value = 5
$ws.Range("C48").Value = 5
$ws.Range("C49").Value = 4
$ws.Range("C51").Value = 1
$ws.Range("C52").Value = 3
$ws.Range("C55").Value = 10
$ws.Range("C55").Style = "Normal"
$ws.Range("C56").Value = 10
$ws.Range("C56").Style = "Normal"
$ws.Range("C57").Value = 10
$ws.Range("C57").Style = "Normal"
$ws.Range("C58").Value = 10
$ws.Range("C58").Style = "Normal"
$ws.Range("C59").Value = 10
$ws.Range("C59").Style = "Normal"
$ws.Range("C60").Value = 10
$ws.Range("C60").Style = "Normal"
$ws.Range("C63").Value = 11
$ws.Range("C63").Style = "Normal"
$ws.Range("C64").Value = 1
$ws.Range("C64").Style = "Normal"
$ws.Range("C65").Value = 3
$ws.Range("C65").Style = "Normal"
$ws.Range("C66").Value = 6
$ws.Range("C66").Style = "Normal"
$ws.Range("C67").Value = 40
$ws.Range("C67").Style = "Normal"
$ws.Range("C72").Value = 2
$ws.Range("C72").Style = "Normal"
$ws.Range("C73").Value = 1
$ws.Range("C73").Style = "Normal"
$ws.Range("C74").Value = 7
$ws.Range("C74").Style = "Normal"
$ws.Range("C75").Value = 47
$ws.Range("C75").Style = "Normal"
$ws.Range("C114").Value = 5
$ws.Range("C115").Value = 4
$ws.Range("C116").Value = 4
$ws.Range("C117").Value = 47
$ws.Range("C121").Value = 3
$ws.Range("C122").Value = 3
$ws.Range("C123").Value = 3
$ws.Range("C124").Value = 30
$ws.Range("C126").Value = 6
$ws.Range("C127").Value = 17

# Update selection/view state
$ws.Range("C77").Select()
